{"js": "// The \"Abstract Title\" paragraph style is removed entirely (it was an\n// unused custom style), the \"Abstract\" style's space-before is bumped\n// from 100 twips (5pt) up to 300 twips (15pt) to match its space-after,\n// and the unused \"Footnote Block Text\" paragraph style is removed.\n\nconst styles = context.document.getStyles();\n\n// Remove the \"Abstract Title\" paragraph style.\nconst abstractTitleStyle = styles.getByName(\"Abstract Title\");\nabstractTitleStyle.delete();\n\n// \"Abstract\" style: space-before 5pt (100 twips) -> 15pt (300 twips).\nconst abstractStyle = styles.getByName(\"Abstract\");\nabstractStyle.paragraphFormat.spaceBefore = 15;\n\n// Remove the \"Footnote Block Text\" paragraph style.\nconst footnoteBlockTextStyle = styles.getByName(\"Footnote Block Text\");\nfootnoteBlockTextStyle.delete();\n\nawait context.sync();\n", "ps1": "# The \"Abstract Title\" paragraph style is removed entirely (it was an\n# unused custom style), the \"Abstract\" style's space-before is bumped\n# from 100 twips (5pt) up to 300 twips (15pt) to match its space-after,\n# and the unused \"Footnote Block Text\" paragraph style is removed.\n\n$d = $word.ActiveDocument\n\n# Remove the \"Abstract Title\" paragraph style.\n$d.Styles(\"Abstract Title\").Delete()\n\n# \"Abstract\" style: space-before 5pt (100 twips) -> 15pt (300 twips).\n$d.Styles(\"Abstract\").ParagraphFormat.SpaceBefore = 15\n\n# Remove the \"Footnote Block Text\" paragraph style.\n$d.Styles(\"Footnote Block Text\").Delete()\n"}
